# "NewProduct" test case: add a second little table (Name / Sale Price /
# Itemnumber / Description) to the right of the existing login table on
# Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the new header + data cells in reading order (this is also the
# order Excel would intern new strings into sharedStrings.xml).
$ws.Range("F1").Value = "Name"
$ws.Range("F2").Value = "John"

$ws.Range("G1").Value = "Sale Price"
$ws.Range("G2").Value = 12
$ws.Range("G2").NumberFormat = "0.00"

$ws.Range("H1").Value = "Itemnumber"
$ws.Range("H2").Value = "Ac10"

$ws.Range("I1").Value = "Description"
$ws.Range("I2").Value = "Butter"

# Column H/I get an explicit (best-fit-like) width, matching the rest of
# the sheet's already-autofitted columns.
$ws.Columns.Item(8).ColumnWidth = 10.5
$ws.Columns.Item(9).ColumnWidth = 9.5

# View state: scroll the window right a bit and leave G2 selected/active.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G2").Select() | Out-Null
